$wb = $excel.ActiveWorkbook
$ws2 = $wb.Worksheets.Item("Sheet2")
foreach ($ws in $wb.Worksheets) {
    $n = $ws.ChartObjects().Count
    if ($n -gt 0) {
        $co = $ws.ChartObjects(1)
        $chart = $co.Chart
        $s = $chart.SeriesCollection(1)
        $s.Values = $ws2.Range("D2:D23")
        $s.XValues = $ws2.Range("A2:A23")
        Write-Host ("new formula=" + $s.Formula)
    }
}
